# Add a new data row (row 71) to each of the 4 worksheets, mirroring the
# structure/format of the existing rows (time, length/ID/checksum hex
# strings, and their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rowsData = @{
    "FE_LFT_#1" = @{
        A = 45857.49570601852
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x34"
        E = "0xf"
        F = 380
        G = "7.598631275147109e+23"
        H = 308
        I = 15
    }
    "FE_LFT_#2" = @{
        A = 45857.49570601852
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x44"
        E = "0xe"
        F = 400
        G = "5.68432987514711e+23"
        H = 324
        I = 14
    }
    "FE_PLT_#1" = @{
        A = 45857.49570601852
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x62"
        E = "0x3"
        F = 110
        G = "5.68631262647114e+23"
        H = 98
        I = 3
    }
    "FE_PLT_#2" = @{
        A = 45857.49570601852
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x62"
        E = "0x3"
        F = 110
        G = "9.85046333984776e+23"
        H = 98
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $rowsData.ContainsKey($name)) {
        continue
    }
    $data = $rowsData[$name]

    $newRow = 71

    # Column A: timestamp, copy the date/time number format used by the
    # preceding row so it renders the same way.
    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.Value = $data.A
    $cellA.NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = [double]$data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
